$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.379369
$ws.Range("H2").Value = 31.138107
$ws.Range("I2").Value = 0.01614698522449884
$ws.Range("J2").Value = 0.01614698522449883
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.07074633333333
$ws.Range("N2").Value = 102.212239
$ws.Range("O2").Value = 0.5171464495142372
$ws.Range("P2").Value = 0.5171464495142373
$ws.Range("Q2").Value = 353.6328482990636
$ws.Range("R2").Value = 3182.695634691573
$ws.Range("S2").Value = 0.00835035607920842
$ws.Range("T2").Value = 0.00835035607920842

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.379369
$ws.Range("H3").Value = 31.138107
$ws.Range("I3").Value = 0.01614698522449884
$ws.Range("J3").Value = 0.01614698522449883
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.685497
$ws.Range("N3").Value = 83.056491
$ws.Range("O3").Value = 0.420227262899125
$ws.Range("P3").Value = 0.4202272628991251
$ws.Range("Q3").Value = 287.357989311393
$ws.Range("R3").Value = 2586.221903802537
$ws.Range("S3").Value = 0.006785403404963759
$ws.Range("T3").Value = 0.006785403404963759

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.379369
$ws.Range("H4").Value = 31.138107
$ws.Range("I4").Value = 0.01614698522449884
$ws.Range("J4").Value = 0.01614698522449883
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.125957666666666
$ws.Range("N4").Value = 12.377873
$ws.Range("O4").Value = 0.06262628758663766
$ws.Range("P4").Value = 0.06262628758663766
$ws.Range("Q4").Value = 42.82483710071233
$ws.Range("R4").Value = 385.423533906411
$ws.Range("S4").Value = 0.001011225740326653
$ws.Range("T4").Value = 0.001011225740326653

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 604.0312093333333
$ws.Range("H5").Value = 1812.093628
$ws.Range("I5").Value = 0.9396797639857967
$ws.Range("J5").Value = 0.9396797639857967
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.07074633333333
$ws.Range("N5").Value = 102.212239
$ws.Range("O5").Value = 0.5171464495142372
$ws.Range("P5").Value = 0.5171464495142373
$ws.Range("Q5").Value = 20579.79411061256
$ws.Range("R5").Value = 185218.1469955131
$ws.Range("S5").Value = 0.4859520536256311
$ws.Range("T5").Value = 0.4859520536256312

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 604.0312093333333
$ws.Range("H6").Value = 1812.093628
$ws.Range("I6").Value = 0.9396797639857967
$ws.Range("J6").Value = 0.9396797639857967
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.685497
$ws.Range("N6").Value = 83.056491
$ws.Range("O6").Value = 0.420227262899125
$ws.Range("P6").Value = 0.4202272628991251
$ws.Range("Q6").Value = 16722.90423390437
$ws.Range("R6").Value = 150506.1381051393
$ws.Range("S6").Value = 0.3948790552214472
$ws.Range("T6").Value = 0.3948790552214472

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 604.0312093333333
$ws.Range("H7").Value = 1812.093628
$ws.Range("I7").Value = 0.9396797639857967
$ws.Range("J7").Value = 0.9396797639857967
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.125957666666666
$ws.Range("N7").Value = 12.377873
$ws.Range("O7").Value = 0.06262628758663766
$ws.Range("P7").Value = 0.06262628758663766
$ws.Range("Q7").Value = 2492.207199054805
$ws.Range("R7").Value = 22429.86479149324
$ws.Range("S7").Value = 0.0588486551387183
$ws.Range("T7").Value = 0.0588486551387183

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.39480333333333
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04417325078970442
$ws.Range("J8").Value = 0.04417325078970442
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.07074633333333
$ws.Range("N8").Value = 102.212239
$ws.Range("O8").Value = 0.5171464495142372
$ws.Range("P8").Value = 0.5171464495142373
$ws.Range("Q8").Value = 967.4321415548875
$ws.Range("R8").Value = 8706.88927399399
$ws.Range("S8").Value = 0.02284403980939761
$ws.Range("T8").Value = 0.02284403980939762

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.39480333333333
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04417325078970442
$ws.Range("J9").Value = 0.04417325078970442
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.685497
$ws.Range("N9").Value = 83.056491
$ws.Range("O9").Value = 0.420227262899125
$ws.Range("P9").Value = 0.4202272628991251
$ws.Range("Q9").Value = 786.12424250059
$ws.Range("R9").Value = 7075.118182505309
$ws.Range("S9").Value = 0.0185628042727141
$ws.Range("T9").Value = 0.0185628042727141

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.39480333333333
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04417325078970442
$ws.Range("J10").Value = 0.04417325078970442
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.125957666666666
$ws.Range("N10").Value = 12.377873
$ws.Range("O10").Value = 0.06262628758663766
$ws.Range("P10").Value = 0.06262628758663766
$ws.Range("Q10").Value = 117.1557565066589
$ws.Range("R10").Value = 1054.40180855993
$ws.Range("S10").Value = 0.002766406707592698
$ws.Range("T10").Value = 0.002766406707592698

